$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12198
$ws1.Range("F4").Value = 54
$ws1.Range("F8").Value = 12119
$ws1.Range("F10").Value = 1190
$ws1.Range("F12").Value = 609
$ws1.Range("F13").Value = 2807
$ws1.Range("F14").Value = 5962
$ws1.Range("F16").Value = 3568

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 12198
$ws4.Range("F5").Value = 54
$ws4.Range("F10").Value = 12119
$ws4.Range("F12").Value = 1190
$ws4.Range("F14").Value = 609
$ws4.Range("F15").Value = 2807
$ws4.Range("F17").Value = 5962
$ws4.Range("F19").Value = 3568

$wb.Save()
